$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Add the missing Finish Date for "Mindset" (row 11)
# Copy the date-formatted cell C11 into D11 first so D11 inherits the
# same number-format style, then overwrite its value.
$ws.Range("C11").Copy($ws.Range("D11"))
$ws.Range("D11").Value = 43847

# Add a new row 12 for "How We Learn"
$ws.Range("A12").Value = "How We Learn "
$ws.Range("B12").Value = "Benedict Carey"

$ws.Range("C11").Copy($ws.Range("C12"))
$ws.Range("C12").Value = 43846

$ws.Range("C11").Copy($ws.Range("D12"))
$ws.Range("D12").Value = 43847

$ws.Range("E12").Value = "learning;psychology;science;neurology"
$ws.Range("F12").Value = "Audio"
$ws.Range("G12").Value = "7 Hrs 21 Mins"

$ws.Range("A13").Select()
